# Add team record (Wins/Losses/Ties) columns to the data sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - reuse the exact formatting of the existing header cells
# (bold, bordered, centered) by copying format from an existing header cell.
$ws.Range("AA1").Copy()
$ws.Range("AC1:AE1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("AC1").Value = "Wins"
$ws.Range("AD1").Value = "Losses"
$ws.Range("AE1").Value = "Ties"

# Fill data rows 2-48 with the team's win/loss/tie record
$lastRow = 48
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 29).Value = 73   # AC - Wins
    $ws.Cells.Item($r, 30).Value = 89   # AD - Losses
    $ws.Cells.Item($r, 31).Value = 0    # AE - Ties
}
